# Exemplo PHP 08/10/2024 2024.2
# Applies two small code-sample corrections:
#   1. Slide 16: "let idade = 16, eleitor;" -> "let idade = 16;"
#   2. Slide 6:  fixes the "carro" object literal declaration and adds
#      ".ano" to the console.log call.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 16 ("Exemplo (Ternário)"): "	let idade = 16, eleitor;"
#                                -> "	let idade = 16;"
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$shape16 = $s16.Shapes.Item(2)

$tr = $shape16.TextFrame.TextRange
$para5 = $tr.Paragraphs(5, 1)
$p5start = $para5.Start
$p5len = $para5.Length

# " idade = 16, eleitor;" is 21 characters, sitting right before the
# paragraph's trailing mark.
$tailStart = $p5start + $p5len - 1 - 21
$tail = $tr.Characters($tailStart, 21)
$tail.Text = " idade = 16;"

# Re-fetch (the text length changed) and split " idade = 16;" into two
# runs: " idade " and "= 16;" (matching the target run layout).
$tr = $shape16.TextFrame.TextRange
$firstPart = $tr.Characters($tailStart, 7)
$secondPart = $tr.Characters($tailStart + 7, 5)
$secondPart.Font.Size = 24

# ---------------------------------------------------------------------
# Slide 6 ("JS – Objetos {}"): fix the "carro" declaration and the
# console.log call.
#   "	const = carro{ano:2001, marca:"volkswagem", "
#     -> "	const carro = { ano:2001, marca:"volkswagem", "
#   "	console.log(carro); "  -> "	console.log(carro.ano); "
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shape6 = $s6.Shapes.Item(2)

$tr = $shape6.TextFrame.TextRange
$para2 = $tr.Paragraphs(2, 1)
$p2start = $para2.Start

# " = " (offset 6, length 3) -> " "
$eqRange = $tr.Characters($p2start + 6, 3)
$eqRange.Text = " "

# "carro" (now at offset 7, length 5) -> "carro = "
$tr = $shape6.TextFrame.TextRange
$para2 = $tr.Paragraphs(2, 1)
$p2start = $para2.Start
$carroRange = $tr.Characters($p2start + 7, 5)
$carroRange.Text = "carro = "

# "{ano:2001, marca:"" (now at offset 15, length 18) -> "{ ano:2001, marca:""
$tr = $shape6.TextFrame.TextRange
$para2 = $tr.Paragraphs(2, 1)
$p2start = $para2.Start
$braceRange = $tr.Characters($p2start + 15, 18)
$braceRange.Text = "{ ano:2001, marca:`""

# Paragraph 7: "	console.log(carro); " -> "	console.log(carro.ano); "
$tr = $shape6.TextFrame.TextRange
$para7 = $tr.Paragraphs(7, 1)
$p7start = $para7.Start
$carro7 = $tr.Characters($p7start + 13, 5)
$carro7.Text = "carro.ano"
